# Daily attendance processing - 2026-01-24 07:56:28
# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in the "Recorded By" column (column G) throughout the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
